$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy row 52 formatting into new row 53 so column styles (date col A vs others) match the rest of the table
$ws.Range("A52:E52").Copy()
$ws.Range("A53:E53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 2: write the corrected dataset (rows 2-53) for columns A-E
$data = @(
    @(39400, 2007, 4.930115226412357, 2008, $null),
    @(39583, 2008, $null, 2009, $null),
    @(39765, 2008, 1.457587285166628, 2009, $null),
    @(39948, 2009, $null, 2010, $null),
    @(40130, 2009, -0.9140166223623458, 2010, $null),
    @(40310, 2010, -2.034793027571991, 2011, 0.6952816881563351),
    @(40494, 2010, 2.585942866987878, 2011, 4.109775046142405),
    @(40676, 2011, 3.880748577052473, 2012, 2.936333428994109),
    @(40862, 2011, 4.253963781362402, 2012, 2.863367440851095),
    @(41044, 2012, 2.174296999091507, 2013, 1.706732094556851),
    @(41228, 2012, 1.752870900283909, 2013, 1.520397254708405),
    @(41409, 2013, -1.942476814417471, 2014, 1.288975737543607),
    @(41592, 2013, -1.479696720105139, 2014, 2.503951807923066),
    @(41774, 2014, 4.926006686093287, 2015, 2.643411312704802),
    @(41957, 2014, 3.900127535411246, 2015, 1.194160460927884),
    @(42137, 2015, 0.7614971083056643, 2016, 2.348955682567344),
    @(42321, 2015, 0.03947433952959933, 2016, 1.459149667419779),
    @(42503, 2016, 3.004126378695804, 2017, 2.358460376580607),
    @(42689, 2016, 2.192778679161944, 2017, 1.586376095629216),
    @(42867, 2017, 2.384525276921168, 2018, 2.349880264276161),
    @(43053, 2017, 3.40836448860673, 2018, 2.570658574505469),
    @(43145, 2018, 1.512747556564698, 2019, 2.136396488383507),
    @(43235, 2018, 2.677874785158663, 2019, 2.507403033330702),
    @(43326, 2018, 2.474496385100733, 2019, 2.141109033939292),
    @(43418, 2018, 2.799070570134488, 2019, 2.479713128614147),
    @(43510, 2019, 2.786825287537487, 2020, 1.555206765808892),
    @(43600, 2019, 3.746278935679004, 2020, 1.868164064786093),
    @(43691, 2019, 4.066325724929976, 2020, 2.042884966056935),
    @(43783, 2019, 4.195393191694419, 2020, 2.359935293525561),
    @(43875, 2020, 1.899045195796845, 2021, 2.205372086670465),
    @(43966, 2020, 3.278383090085746, 2021, 2.401026764575831),
    @(44068, 2020, 2.133862376612439, 2021, 1.412539570439053),
    @(44159, 2020, 1.666553973046048, 2021, -0.4512719783814068),
    @(44251, 2021, -1.801793552285647, 2022, 1.696367270021448),
    @(44341, 2021, 2.877357105987888, 2022, 2.770626214993133),
    @(44432, 2021, 1.773820722495745, 2022, 2.446369413529137),
    @(44525, 2021, 1.879266440112803, 2022, 1.081814991510499),
    @(44617, 2022, -0.1058999733161259, 2023, 2.729902493405767),
    @(44706, 2022, -0.4001895765463725, 2023, 2.524513388369543),
    @(44798, 2022, -2.404913754290983, 2023, 1.366747064788676),
    @(44890, 2022, -2.620683231370946, 2023, -1.174318230871441),
    @(44981, 2023, -3.321926039826262, 2024, 1.506553781929298),
    @(45071, 2023, -2.321114556364801, 2024, 1.637366184014355),
    @(45163, 2023, -2.901570548279864, 2024, 1.175840608617551),
    @(45254, 2023, -3.036556262700274, 2024, 0.07123445333143685),
    @(45345, 2024, -1.953081240592103, 2025, 0.8938306900690307),
    @(45436, 2024, -0.9245282159112467, 2025, 0.9277865862836965),
    @(45534, 2024, -2.859191689251428, 2025, -0.3916886409131459),
    @(45618, 2024, -2.953443685011514, 2025, -1.196842846539037),
    @(45713, 2025, -0.2732731229103447, 2026, 0.7812052424394755),
    @(45800, 2025, 0.2384226118222088, 2026, 0.8054161303035379),
    @(45891, 2025, -1.131442475565558, 2026, -0.4808467302335195),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($c = 1; $c -le 5; $c++) {
        $val = $rowVals[$c - 1]
        if ($null -eq $val) {
            $ws.Cells.Item($r, $c).ClearContents()
        } else {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

Write-Output "Edit complete"